$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new "16 nm" feature-size columns, one in each metric group ---
# (Insert() copies formatting from the left neighbour by default, which is wrong for
#  these "interior" columns, so re-paste the format from the correct neighbour on the
#  right immediately after each insert.)
$ws.Columns("E:E").Insert()
$ws.Range("F1:F7").Copy()
$ws.Range("E1:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Columns("I:I").Insert()
$ws.Range("J1:J7").Copy()
$ws.Range("I1:I7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new "16 nm" feature-size values in row 3 ---
$ws.Range("E3").Value = 16
$ws.Range("I3").Value = 16

# --- Row 6: nvdla (its area/power baseline only exists starting at 16nm) ---
$ws.Range("D6").Clear()
$ws.Range("H6").Clear()
$ws.Range("C6").Value = "nvdla"
$ws.Range("E6").Value = 1000000
$ws.Range("F6").Formula = "=E6*(F`$3/`$E`$3)^2"
$ws.Range("G6").Formula = "=`$E6*(G`$3/`$E`$3)^2"
$ws.Range("I6").Value = 48
$ws.Range("J6").Formula = "=I6*(J`$3/`$I`$3)^2"
$ws.Range("K6").Formula = "=I6*(K`$3/`$I`$3)^2"

# --- Row 7: sdp (same shape as row 4 gemm) ---
$ws.Range("C7").Value = "sdp"
$ws.Range("D7").Value = 54288
$ws.Range("F7").Formula = "=D7*(F`$3/`$D`$3)^2"
$ws.Range("G7").Formula = "=`$D7*(G`$3/`$D`$3)^2"
$ws.Range("H7").Value = 34
$ws.Range("J7").Formula = "=H7*(J`$3/`$H`$3)^2"
$ws.Range("K7").Formula = "=H7*(K`$3/`$H`$3)^2"

Write-Host "done"
